$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New experiment rows appended to the log (rows 11-18)
# Columns: A = file_name, B = experiment_date (text, quote-prefixed), C = description

# Row 11 uses the same date style as row 10 (B10, style s="3")
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").Value = "growth curves aureus stampersruw fractie 21mrt2022.xlsx"
$ws.Range("B11").Value = "'21-3-2022"
$ws.Range("C11").Value = "S.aureus with red tulip anther extract"

# Row 12 uses the same date style as row 10 (B10, style s="3")
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Value = "growth curves aures huas028 22mrt2022.xlsx"
$ws.Range("B12").Value = "'22-3-2022"
$ws.Range("C12").Value = "S.aureus with red tulip stamen extract"

# Row 13 uses the same date style as row 10 (B10, style s="3")
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value = "growth curves ecoli huas028 4april2022.xlsx"
$ws.Range("B13").Value = "'4-4-2022"
$ws.Range("C13").Value = "E.coli with red tulip stamen extract"

# Row 14 uses the same date style as row 9 (B9, style s="2")
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = "growth curves kpneumoniae huas028 5april2022.xlsx"
$ws.Range("B14").Value = "'5-4-2022"
$ws.Range("C14").Value = "K.pneumoniae with red tulip stamen extract"

# Row 15 uses the same date style as row 9 (B9, style s="2")
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Value = "growth curves paeruginosa huas028 7april2022.xlsx"
$ws.Range("B15").Value = "'7-4-2022"
$ws.Range("C15").Value = "P.aeruginosa with red tulip stamen extract"

# Row 16 uses the same date style as row 9 (B9, style s="2")
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = "growth curves aureus huas033 28april2022.xls"
$ws.Range("B16").Value = "'28-4-2022"
$ws.Range("C16").Value = "S.aureus with Delphinidin-3-rutinoside"

# Row 17 uses the same date style as row 9 (B9, style s="2")
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = "growth curves ecoli huas034 16mei2022.xlsx"
$ws.Range("B17").Value = "'16-5-2022"
$ws.Range("C17").Value = "E.coli with Delphinidin-3-rutinoside"

# Row 18 uses the same date style as row 9 (B9, style s="2")
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = "growth curves kpneumoniae huas034 31mei2022.xlsx"
$ws.Range("B18").Value = "'31-5-2022"
$ws.Range("C18").Value = "K.pneumoniae with Delphinidin-3-rutinoside"

# Restore the active selection to match the author's final cursor position
$ws.Range("D13").Select() | Out-Null
